# Add a new "cfop" column to the "PI hours" sheet, and a new "cfop hours"
# worksheet (placed after "unit(accumulative) hours"), mirroring the
# layout/format of the existing "PI hours" sheet.

$wb = $excel.ActiveWorkbook
$pi = $wb.Worksheets.Item("PI hours")

# --- 1. Add "cfop" column (G) to the "PI hours" sheet -----------------
$pi.Range("G1").Value = "cfop"
$pi.Range("F1").Copy()
$pi.Range("G1").PasteSpecial(-4122)   # xlPasteFormats

$cfopByName = @{
    "Naira Hovakimyan"    = "['cfop_NH']"
    "Seth Hutchinson"     = "['cfop_HUTCHINSON']"
    "Romit Roy Choudhury" = "['cfop_CHOUDHURY', 'cfop_RRC']"
    "Hae-Won Park"        = "['cfop_PARK']"
    "Girish Chowdhary"    = "['cfop_GC']"
    "Michael Selig"       = "['cfop_SELIG']"
}

for ($r = 2; $r -le 7; $r++) {
    $name = $pi.Cells.Item($r, 2).Value()
    $pi.Cells.Item($r, 7).Value = $cfopByName[$name]
}

# --- 2. Add the new "cfop hours" worksheet, after the last sheet ------
$count = $wb.Worksheets.Count()
$lastSheet = $wb.Worksheets.Item($count)
$cfop = $wb.Worksheets.Add($null, $lastSheet)
$cfop.Name = "cfop hours"

# Header row (cfop / hours / percentage), formatted like PI hours' header.
$cfop.Range("B1").Value = "cfop"
$cfop.Range("C1").Value = "hours"
$cfop.Range("D1").Value = "percentage"
$pi.Range("B1:D1").Copy()
$cfop.Range("B1:D1").PasteSpecial(-4122)

$cfopRows = @(
    @("cfop_HUTCHINSON", 56,   30.27027027027027),
    @("cfop_RRC",         49,  26.48648648648649),
    @("cfop_NH",          35,  18.91891891891892),
    @("cfop_CHOUDHURY",   24.5, 13.24324324324324),
    @("cfop_PARK",        13,  7.027027027027027),
    @("cfop_GC",          6,   3.243243243243243),
    @("cfop_SELIG",       1.5, 0.8108108108108109)
)

$r = 2
foreach ($row in $cfopRows) {
    $cfop.Cells.Item($r, 1).Value = $r - 2
    $cfop.Cells.Item($r, 2).Value = $row[0]
    $cfop.Cells.Item($r, 3).Value = $row[1]
    $cfop.Cells.Item($r, 4).Value = $row[2]
    $r++
}

# Index column (A) formatted like PI hours' index column.
$pi.Range("A2").Copy()
$cfop.Range("A2:A8").PasteSpecial(-4122)

# Restore the original active sheet / selection.
$pi.Activate()
